$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for SVR parameters
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# New data values
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Clear leftover empty formatted cell in row 13 (A13) so the row disappears
$ws.Range("A13").Clear()

# Remove the now-duplicate/unused cell style (old index 1, identical to default)
# by clearing it from every cell/column that referenced it, collapsing cellXfs
$ws.Columns("B:C").ClearFormats()
$dupStyleCells = @("B1","C1","H1","H2","A5","A6","A7","A8","A11","D14","E14","F14","G14","I14","D15","E15","F15","G15","I15","D16","E16","F16","G16","I16")
foreach ($c in $dupStyleCells) {
    $ws.Range($c).Style = "Normal"
}

# Update selected cell
$ws.Range("I9").Select()
